# "data update in excel"
# Update the failed-test detail rows so that both TC1 and TC2 report the
# same ERR-5005 "Internal server error" failure, widen column C, and
# leave the selection on the freshly edited range (A2:XFD7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$errorCode    = "ERR-5005"
$errorMessage = "Internal server error"
$logEntry     = "ERROR - ERR-5005: Internal server error in module Hatity. User: kasie.ankunding, SessionID: 0f210e9c-8293-4087-bfed-98d284d7bc54"

# Row 2 (TC1)
$ws.Range("B2").Value = $errorCode
$ws.Range("C2").Value = $errorMessage
$ws.Range("D2").Value = $logEntry

# Row 3 (TC2)
$ws.Range("B3").Value = $errorCode
$ws.Range("C3").Value = $errorMessage
$ws.Range("D3").Value = $logEntry

# Widen column C (errorMessage) to fit the new content.
$ws.Columns.Item(3).ColumnWidth = 29.5

# Leave the selection on the updated rows, as in the saved workbook.
$ws.Range("A2:XFD7").Select() | Out-Null
